# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Wed May 29 02:47:01 UTC 2024 with GitHub Actions"
# Source data is scraped text (coinranking.com), so every cell on the sheet
# is stored as text, even the numeric-looking price figures. Cells whose new
# value would otherwise be auto-coerced to a number by Excel are pre-
# formatted as Text ("@") so the literal digits (e.g. trailing zeros,
# leading "0.000...") survive exactly as scraped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.601.49'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '3.846.14'
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.92'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.09'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = '3.848.01'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000277'
$ws.Range("E13").Value = '  +10.49%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.85'
$ws.Range("E14").Value = '  -2.11%  '
$ws.Range("D15").Value = '4.493.21'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").Value = '3.848.22'
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("D17").Value = '68.645.23'
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.33'
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.33'
$ws.Range("E19").Value = '  -3.08%  '
$ws.Range("E20").Value = '  -0.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.93'
$ws.Range("E21").Value = '  +1.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '472.01'
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.724'
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.29'
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.40'
$ws.Range("E28").Value = '  +3.21%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").Value = '3.998.70'
$ws.Range("E31").Value = '  -0.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.70'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.32'
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("E34").Value = '  -1.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.26'
$ws.Range("E35").Value = '  -2.32%  '
$ws.Range("D36").Value = '3.813.67'
$ws.Range("E36").Value = '  -0.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +14.97%  '
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.140'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.02'
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.90'
$ws.Range("E41").Value = '  -1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.315'
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  -1.35%  '
$ws.Range("B45").Value = 'FLOKI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000300'
$ws.Range("E45").Value = '  +10.14%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '419.78'
$ws.Range("E47").Value = '  -2.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.65'
$ws.Range("E48").Value = '  +0.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '46.89'
$ws.Range("E49").Value = '  -2.00%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.67'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.02'
$ws.Range("E51").Value = '  +4.61%  '
